$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 89, pushing the existing rows 89-115 down to 90-116.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly price observation
# (same market/category/quality metadata as the surrounding rows, new date
# and volume, and the same min/max/avg price + $/Kg figures as before).
$ws.Range("A89").Value = 8
$ws.Range("B89").Value = "Terminal La Palmera de La Serena"
$ws.Range("C89").Value = "Coquimbo"
$ws.Range("D89").Value = 44468
$ws.Range("E89").Value = 4
$ws.Range("F89").Value = 100112037
$ws.Range("G89").Value = "Cebollín"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 3160
$ws.Range("K89").Value = 900
$ws.Range("L89").Value = 1000
$ws.Range("M89").Value = 950
$ws.Range("N89").Value = "`$/paquete 6 unidades"
$ws.Range("O89").Value = "Provincia del Elquí"
$ws.Range("P89").Value = 158
$ws.Range("Q89").Value = 6
$ws.Range("R89").Value = "Hortaliza"
